# Updates odds values on Sheet1 as per the FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("J3").Value = 2.38
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 7
$ws.Range("AE3").Value = 19
$ws.Range("AM3").Value = 451
$ws.Range("AN3").Value = 3.5
$ws.Range("AO3").Value = 9

# Row 6
$ws.Range("J6").Value = 3.1
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("BB6").Value = 351

# Row 7
$ws.Range("O7").Value = 1.62
$ws.Range("P7").Value = 2.2

# Row 8
$ws.Range("G8").Value = 1.55
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 5.75
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 2
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.83
$ws.Range("X8").Value = 7.5
$ws.Range("AG8").Value = 15
$ws.Range("AH8").Value = 29
$ws.Range("AL8").Value = 41
$ws.Range("AO8").Value = 8
$ws.Range("AU8").Value = 8.5
$ws.Range("AW8").Value = 7
$ws.Range("AX8").Value = 29
$ws.Range("AY8").Value = 34
$ws.Range("AZ8").Value = 101
$ws.Range("BA8").Value = 126
$ws.Range("BB8").Value = 251

# Row 9
$ws.Range("G9").Value = 2.63
$ws.Range("I9").Value = 2.7
$ws.Range("J9").Value = 3.4
$ws.Range("L9").Value = 3.5
$ws.Range("Y9").Value = 10
$ws.Range("Z9").Value = 26
$ws.Range("AI9").Value = 11
$ws.Range("AJ9").Value = 29
$ws.Range("AN9").Value = 4.5
$ws.Range("AO9").Value = 15
$ws.Range("AW9").Value = 4.75
$ws.Range("AX9").Value = 17
